$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

# Clear the placeholder quote-strings from cells B5, C6, B7, C7
$ws.Range("B5").ClearContents()
$ws.Range("C6").ClearContents()
$ws.Range("B7").ClearContents()
$ws.Range("C7").ClearContents()

# Update the active selection to B7
$ws.Range("B7").Select()
